# Cambios Año y Catalogos
# Applies budget (M:T) / actual (U:AB) column updates for rows 2, 4, 45
# and targeted value fixes for rows 14, 49, plus the A-column
# renumbering (Numero de orden / catalog renumbering) for rows 11-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: add Budget (M:T) values, and make Actual (U:AB) mirror them ---
$ws.Cells.Item(2, 13).Value = 13391.57   # M2  Servicio_budget
$ws.Cells.Item(2, 14).Value = 279350.01  # N2  Equipo_budget
$ws.Cells.Item(2, 15).Value = 0          # O2  Protectores de Cable_budget
$ws.Cells.Item(2, 16).Value = 0          # P2  Capilar_budget
$ws.Cells.Item(2, 17).Value = 0          # Q2  Equipo Superficie_budget
$ws.Cells.Item(2, 18).Value = 0          # R2  Desarenador_budget
$ws.Cells.Item(2, 19).Value = 63774.22   # S2  Cable Nuevo_budget
$ws.Cells.Item(2, 20).Value = 0          # T2  B&H_budget

$ws.Cells.Item(2, 21).Value = 13391.57   # U2  Servicio_Actual
$ws.Cells.Item(2, 22).Value = 279350.01  # V2  Equipo_Actual
$ws.Cells.Item(2, 23).Value = 0          # W2  Protectores de Cable_Actual
$ws.Cells.Item(2, 27).Value = 63774.22   # AA2 Cable Nuevo_Actual
$ws.Cells.Item(2, 28).Value = 0          # AB2 B&H_Actual

# --- Row 4: add Budget (M:T) values, and make Actual (U:AB) mirror them ---
$ws.Cells.Item(4, 13).Value = 33666.72   # M4  Servicio_budget
$ws.Cells.Item(4, 14).Value = 341703.8   # N4  Equipo_budget
$ws.Cells.Item(4, 15).Value = 0          # O4  Protectores de Cable_budget
$ws.Cells.Item(4, 16).Value = 0          # P4  Capilar_budget
$ws.Cells.Item(4, 17).Value = 0          # Q4  Equipo Superficie_budget
$ws.Cells.Item(4, 18).Value = 0          # R4  Desarenador_budget
$ws.Cells.Item(4, 19).Value = 0          # S4  Cable Nuevo_budget
$ws.Cells.Item(4, 20).Value = 0          # T4  B&H_budget

$ws.Cells.Item(4, 21).Value = 33666.72   # U4  Servicio_Actual
$ws.Cells.Item(4, 22).Value = 341703.8   # V4  Equipo_Actual
$ws.Cells.Item(4, 23).Value = 0          # W4  Protectores de Cable_Actual
$ws.Cells.Item(4, 27).Value = 0          # AA4 Cable Nuevo_Actual
$ws.Cells.Item(4, 28).Value = 0          # AB4 B&H_Actual

# --- A-column renumbering (catalog numbers) for rows 11-49 ---
$newA = @{
    11 = 11; 12 = 12; 13 = 14; 14 = 15; 15 = 16; 16 = 18; 17 = 19; 18 = 20;
    19 = 21; 20 = 22; 21 = 23; 22 = 24; 23 = 25; 24 = 26; 25 = 27; 26 = 28;
    27 = 29; 28 = 30; 29 = 31; 30 = 32; 31 = 33; 32 = 34; 33 = 35; 34 = 36;
    35 = 37; 36 = 38; 37 = 41; 38 = 42; 39 = 43; 40 = 44; 41 = 45; 42 = 46;
    43 = 47; 44 = 48; 45 = 49; 46 = 50; 47 = 51; 48 = 52; 49 = 53
}
foreach ($r in $newA.Keys) {
    $ws.Cells.Item($r, 1).Value = $newA[$r]
}

# --- Row 14: Equipo_budget (N14) / Equipo_Actual (V14) corrected ---
$ws.Cells.Item(14, 14).Value = 334593.33  # N14
$ws.Cells.Item(14, 22).Value = 334593.33  # V14

# --- Row 45: B&H_cotizacion (L45) updated + Budget (M:T) added, Actual mirrors Budget ---
$ws.Cells.Item(45, 12).Value = 74021.34   # L45 B&H_cotizacion

$ws.Cells.Item(45, 13).Value = 26512.21   # M45 Servicio_budget
$ws.Cells.Item(45, 14).Value = 337917.65  # N45 Equipo_budget
$ws.Cells.Item(45, 15).Value = 60888.94   # O45 Protectores de Cable_budget
$ws.Cells.Item(45, 16).Value = 0          # P45 Capilar_budget
$ws.Cells.Item(45, 17).Value = 0          # Q45 Equipo Superficie_budget
$ws.Cells.Item(45, 18).Value = 0          # R45 Desarenador_budget
$ws.Cells.Item(45, 19).Value = 0          # S45 Cable Nuevo_budget
$ws.Cells.Item(45, 20).Value = 74021.34   # T45 B&H_budget

$ws.Cells.Item(45, 21).Value = 26512.21   # U45 Servicio_Actual
$ws.Cells.Item(45, 22).Value = 337917.65  # V45 Equipo_Actual
$ws.Cells.Item(45, 23).Value = 60888.94   # W45 Protectores de Cable_Actual
$ws.Cells.Item(45, 24).Value = 0          # X45 Capilar_Actual
$ws.Cells.Item(45, 28).Value = 74021.34   # AB45 B&H_Actual

# --- Row 49: Servicio_budget (M49) / Servicio_Actual (U49) corrected ---
$ws.Cells.Item(49, 13).Value = 136868.28  # M49
$ws.Cells.Item(49, 21).Value = 136868.28  # U49
